$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helper: write literal text into a range without Excel's automatic
# type inference turning boolean-looking ("false") or date-looking
# ("2019-12-01") text into a native Boolean/Date. We build the text via a
# formula (a quoted string literal always stays text) then flatten the
# formula down to a static value with a values-only paste, which keeps the
# shared-string text type and does not touch cell formatting/styles.
function Set-LiteralText {
    param($range, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# ---- Existing rows 2-4 ----
# Elapsed Time (D) changes; these are plain numbers so a normal .Value
# assignment is safe.
$ws.Range("D2").Value = 6
$ws.Range("D3").Value = 10
$ws.Range("D4").Value = 7

# Date (E) text changes from "2019-10-21" to "2019-11-04" for all three
# rows at once, so the shared string is updated in place instead of
# leaving the old string orphaned.
Set-LiteralText $ws.Range("E2:E4") "2019-11-04"

# ---- New rows 5-9 ----
$ws.Range("A5").Value = "chaotsai@stonybrook.edu"
$ws.Range("B5").Value = 1
Set-LiteralText $ws.Range("C5") "false"
$ws.Range("D5").Value = 32
Set-LiteralText $ws.Range("E5") "2019-12-01"

$ws.Range("A6").Value = "chaotsai@stonybrook.edu"
$ws.Range("B6").Value = 2
Set-LiteralText $ws.Range("C6") "false"
$ws.Range("D6").Value = 55
Set-LiteralText $ws.Range("E6") "2019-12-01"

$ws.Range("A7").Value = "ikleiman@stonybrook.edu"
$ws.Range("B7").Value = 2
Set-LiteralText $ws.Range("C7") "false"
$ws.Range("D7").Value = 8
Set-LiteralText $ws.Range("E7") "2019-12-03"

$ws.Range("A8").Value = "vlgarcia@stonybrook.edu"
$ws.Range("B8").Value = 1
Set-LiteralText $ws.Range("C8") "false"
$ws.Range("D8").Value = 7
Set-LiteralText $ws.Range("E8") "2019-12-03"

$ws.Range("A9").Value = "vlgarcia@stonybrook.edu"
$ws.Range("B9").Value = 2
Set-LiteralText $ws.Range("C9") "false"
$ws.Range("D9").Value = 10
Set-LiteralText $ws.Range("E9") "2019-12-03"
